$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.711.68"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "3.517.23"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'623.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.03%  "
$ws.Range("D6").Value = "'171.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").Value = "3.510.04"
$ws.Range("E7").Value = "  -1.49%  "
$ws.Range("D8").Value = "'0.608"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "'0.200"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("E11").Value = "  -2.77%  "
$ws.Range("D12").Value = "'0.585"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").Value = "'46.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").Value = "'0.0000276"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "4.076.95"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").Value = "'8.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "'605.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").Value = "3.510.70"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").Value = "70.748.67"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("D21").Value = "'17.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("D22").Value = "'0.879"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").Value = "'9.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").Value = "'97.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").Value = "'15.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.06%  "
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "'2.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.68%  "
$ws.Range("D29").Value = "'33.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("D30").Value = "'9.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("D31").Value = "'3.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.45%  "
$ws.Range("D32").Value = "'8.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("D34").Value = "'6.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.96%  "
$ws.Range("D35").Value = "'621.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.40%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.0491"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.71%  "
$ws.Range("B37").Value = "Cosmos"
$ws.Range("C37").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D37").Value = "'10.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "'0.0993"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'56.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'3.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.22%  "
$ws.Range("D42").Value = "'0.142"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("D43").Value = "3.332.19"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").Value = "0.0₃0723"
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").Value = "'0.310"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.61%  "
$ws.Range("D47").Value = "'31.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("D48").Value = "'2.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.43%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "'134.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.50%  "
